$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing values (but keep cell formatting/styles) so the shared-string
# table gets rebuilt fresh, in the exact order we write new values below.
$ws.UsedRange.ClearContents()

# Header row
$ws.Range("B1").Value = "feature"
$ws.Range("C1").Value = "impDecrease"
$ws.Range("D1").Value = "Model"

# Column B (feature names) top-to-bottom for all data rows
$ws.Range("B2").Value = "sstk_std"
$ws.Range("B3").Value = "teq"
$ws.Range("B4").Value = "seq"
$ws.Range("B5").Value = "rest_sum_diff"
$ws.Range("B6").Value = "xopr"
$ws.Range("B7").Value = "ceq"
$ws.Range("B8").Value = "caps"
$ws.Range("B9").Value = "tstk"
$ws.Range("B10").Value = "fopo_std"
$ws.Range("B11").Value = "sec_trt1m_std"
$ws.Range("B12").Value = "ceqt"
$ws.Range("B13").Value = "icapt"
$ws.Range("B14").Value = "rest_count_of_diffs"
$ws.Range("B15").Value = "lse"
$ws.Range("B16").Value = "pi_std"
$ws.Range("B17").Value = "revt"
$ws.Range("B18").Value = "sstk"
$ws.Range("B19").Value = "st_per_growth"
$ws.Range("B20").Value = "at"
$ws.Range("B21").Value = "rest_count"
$ws.Range("B22").Value = "cogs"
$ws.Range("B23").Value = "gp"
$ws.Range("B24").Value = "xsga"
$ws.Range("B25").Value = "rect_std"
$ws.Range("B26").Value = "cogs_std"
$ws.Range("B27").Value = "lct"
$ws.Range("B28").Value = "dilavx_std"
$ws.Range("B29").Value = "tstk_std"
$ws.Range("B30").Value = "auop"
$ws.Range("B31").Value = "dlc"
$ws.Range("B32").Value = "xsga_std"
$ws.Range("B33").Value = "at"
$ws.Range("B34").Value = "icapt_std"
$ws.Range("B35").Value = "dltt_std"
$ws.Range("B36").Value = "ceqt"
$ws.Range("B37").Value = "pi_std"
$ws.Range("B38").Value = "st_per_growth"
$ws.Range("B39").Value = "dltr"
$ws.Range("B40").Value = "sstk_std"
$ws.Range("B41").Value = "caps"
$ws.Range("B42").Value = "np_std"
$ws.Range("B43").Value = "st_per_currentToMax"
$ws.Range("B44").Value = "cogs"
$ws.Range("B45").Value = "siv_std"
$ws.Range("B46").Value = "gvkey"
$ws.Range("B47").Value = "invch_std"
$ws.Range("B48").Value = "ivncf"
$ws.Range("B49").Value = "spce"
$ws.Range("B50").Value = "dvt_std"
$ws.Range("B51").Value = "sec_trt1m_std"
$ws.Range("B52").Value = "re"
$ws.Range("B53").Value = "fopo_std"
$ws.Range("B54").Value = "xsga"
$ws.Range("B55").Value = "lt_std"
$ws.Range("B56").Value = "rat_spcsrc"

# Column D (model name) top-to-bottom for all data rows
$ws.Range("D2").Value = "Random Forest"
$ws.Range("D3").Value = "Random Forest"
$ws.Range("D4").Value = "Random Forest"
$ws.Range("D5").Value = "Random Forest"
$ws.Range("D6").Value = "Random Forest"
$ws.Range("D7").Value = "Random Forest"
$ws.Range("D8").Value = "Random Forest"
$ws.Range("D9").Value = "Random Forest"
$ws.Range("D10").Value = "Random Forest"
$ws.Range("D11").Value = "Random Forest"
$ws.Range("D12").Value = "Random Forest"
$ws.Range("D13").Value = "Random Forest"
$ws.Range("D14").Value = "Random Forest"
$ws.Range("D15").Value = "Random Forest"
$ws.Range("D16").Value = "Random Forest"
$ws.Range("D17").Value = "Random Forest"
$ws.Range("D18").Value = "Random Forest"
$ws.Range("D19").Value = "Random Forest"
$ws.Range("D20").Value = "Random Forest"
$ws.Range("D21").Value = "Random Forest"
$ws.Range("D22").Value = "Random Forest"
$ws.Range("D23").Value = "Random Forest"
$ws.Range("D24").Value = "Random Forest"
$ws.Range("D25").Value = "Random Forest"
$ws.Range("D26").Value = "Random Forest"
$ws.Range("D27").Value = "Random Forest"
$ws.Range("D28").Value = "Random Forest"
$ws.Range("D29").Value = "Random Forest"
$ws.Range("D30").Value = "Random Forest"
$ws.Range("D31").Value = "Random Forest"
$ws.Range("D32").Value = "Random Forest"
$ws.Range("D33").Value = "Gradient Boosting"
$ws.Range("D34").Value = "Gradient Boosting"
$ws.Range("D35").Value = "Gradient Boosting"
$ws.Range("D36").Value = "Gradient Boosting"
$ws.Range("D37").Value = "Gradient Boosting"
$ws.Range("D38").Value = "Gradient Boosting"
$ws.Range("D39").Value = "Gradient Boosting"
$ws.Range("D40").Value = "Gradient Boosting"
$ws.Range("D41").Value = "Gradient Boosting"
$ws.Range("D42").Value = "Gradient Boosting"
$ws.Range("D43").Value = "Gradient Boosting"
$ws.Range("D44").Value = "Gradient Boosting"
$ws.Range("D45").Value = "Gradient Boosting"
$ws.Range("D46").Value = "Gradient Boosting"
$ws.Range("D47").Value = "Gradient Boosting"
$ws.Range("D48").Value = "Gradient Boosting"
$ws.Range("D49").Value = "Gradient Boosting"
$ws.Range("D50").Value = "Gradient Boosting"
$ws.Range("D51").Value = "Gradient Boosting"
$ws.Range("D52").Value = "Gradient Boosting"
$ws.Range("D53").Value = "Gradient Boosting"
$ws.Range("D54").Value = "Gradient Boosting"
$ws.Range("D55").Value = "Gradient Boosting"
$ws.Range("D56").Value = "Manual Addition"

# Column C (importance values)
$ws.Range("C2").Value = 0.05152090202630415
$ws.Range("C3").Value = 0.04210914267773873
$ws.Range("C4").Value = 0.04113480038972558
$ws.Range("C5").Value = 0.02964086773085264
$ws.Range("C6").Value = 0.02894210148331851
$ws.Range("C7").Value = 0.02673327775011868
$ws.Range("C8").Value = 0.02374550671631139
$ws.Range("C9").Value = 0.02261428377865916
$ws.Range("C10").Value = 0.02202374152521238
$ws.Range("C11").Value = 0.02085166940573271
$ws.Range("C12").Value = 0.02081067191182673
$ws.Range("C13").Value = 0.02018532661059256
$ws.Range("C14").Value = 0.02006966905846362
$ws.Range("C15").Value = 0.01844682961883291
$ws.Range("C16").Value = 0.0179690776523493
$ws.Range("C17").Value = 0.01768615437752864
$ws.Range("C18").Value = 0.01747198953836201
$ws.Range("C19").Value = 0.01659163279284432
$ws.Range("C20").Value = 0.01653665368324137
$ws.Range("C21").Value = 0.01653338496751571
$ws.Range("C22").Value = 0.01642483665464855
$ws.Range("C23").Value = 0.01639626005398251
$ws.Range("C24").Value = 0.01488677887659137
$ws.Range("C25").Value = 0.01412363442861623
$ws.Range("C26").Value = 0.01403655662386116
$ws.Range("C27").Value = 0.01353114676459853
$ws.Range("C28").Value = 0.01110092391154593
$ws.Range("C29").Value = 0.01082068248664309
$ws.Range("C30").Value = 0.01058660278150156
$ws.Range("C31").Value = 0.01045068502105618
$ws.Range("C32").Value = 0.01028382202854827
$ws.Range("C33").Value = 0.2811118858614863
$ws.Range("C34").Value = 0.05779940637772469
$ws.Range("C35").Value = 0.05709844912089636
$ws.Range("C36").Value = 0.05601291370908551
$ws.Range("C37").Value = 0.04573626607771975
$ws.Range("C38").Value = 0.03871699840744534
$ws.Range("C39").Value = 0.03449255688705752
$ws.Range("C40").Value = 0.02989692735997116
$ws.Range("C41").Value = 0.0288112328024997
$ws.Range("C42").Value = 0.02489605040799534
$ws.Range("C43").Value = 0.02419734966459052
$ws.Range("C44").Value = 0.0238147901688538
$ws.Range("C45").Value = 0.02220252823938738
$ws.Range("C46").Value = 0.02171131822313106
$ws.Range("C47").Value = 0.02033590110204956
$ws.Range("C48").Value = 0.02014301674406948
$ws.Range("C49").Value = 0.01784287285017829
$ws.Range("C50").Value = 0.01739649606700123
$ws.Range("C51").Value = 0.01462492956957438
$ws.Range("C52").Value = 0.01279937493963577
$ws.Range("C53").Value = 0.01141585968354888
$ws.Range("C54").Value = 0.01123499085308485
$ws.Range("C55").Value = 0.01120823540808618
$ws.Range("C56").Value = 0.001111

# Column A (rank index)
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18
$ws.Range("A21").Value = 19
$ws.Range("A22").Value = 20
$ws.Range("A23").Value = 21
$ws.Range("A24").Value = 22
$ws.Range("A25").Value = 23
$ws.Range("A26").Value = 24
$ws.Range("A27").Value = 25
$ws.Range("A28").Value = 26
$ws.Range("A29").Value = 27
$ws.Range("A30").Value = 28
$ws.Range("A31").Value = 29
$ws.Range("A32").Value = 30
$ws.Range("A33").Value = 31
$ws.Range("A34").Value = 32
$ws.Range("A35").Value = 33
$ws.Range("A36").Value = 34
$ws.Range("A37").Value = 35
$ws.Range("A38").Value = 36
$ws.Range("A39").Value = 37
$ws.Range("A40").Value = 38
$ws.Range("A41").Value = 39
$ws.Range("A42").Value = 40
$ws.Range("A43").Value = 41
$ws.Range("A44").Value = 42
$ws.Range("A45").Value = 43
$ws.Range("A46").Value = 44
$ws.Range("A47").Value = 45
$ws.Range("A48").Value = 46
$ws.Range("A49").Value = 47
$ws.Range("A50").Value = 48
$ws.Range("A51").Value = 49
$ws.Range("A52").Value = 50
$ws.Range("A53").Value = 51
$ws.Range("A54").Value = 52
$ws.Range("A55").Value = 53
$ws.Range("A56").Value = 54

# Extend the index-column style to the newly added rows
$ws.Range("A2").Copy()
$ws.Range("A50:A56").PasteSpecial(-4122)
